$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44931
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101006
$ws.Cells.Item($row, 10).Value = "Breva"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 50
$ws.Cells.Item($row, 14).Value = 18000
$ws.Cells.Item($row, 15).Value = 18000
$ws.Cells.Item($row, 16).Value = 18000
$ws.Cells.Item($row, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item($row, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item($row, 19).Value = 3600
$ws.Cells.Item($row, 20).Value = 5
